# "Add files via upload" — adds two new data rows (26 and 27) to the
# "2025LCK春季常规赛" sheet, following the same pattern as the existing
# blocks of rows (A/B = region/team, G = condition text, K = stake size,
# L = remark), and leaves the selection on the last cell touched (L27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 26 --------------------------------------------------------
$ws.Range("A26").Value = "NACL"
$ws.Range("B26").Value = "CCG,WINTER"
$ws.Range("G26").Value = ">32"
$ws.Range("K26").Value = "600"
$ws.Range("L26").Value = "只针对第一把"

# ---- Row 27 ----------------------------------------------------------
$ws.Range("B27").Value = "near,apex"
$ws.Range("C27").Value = 10
$ws.Range("G27").Value = ">29"
$ws.Range("I27").Value = 300
$ws.Range("K27").Value = "800"
$ws.Range("L27").Value = "只针对第一二把，高楼有就i对冲"

# Leave the selection where the author's last edit was.
$ws.Range("L27").Select()
